# Update "想去人数" (F column) values on the sheets that list individual
# event entries ("展览" and "全部类型"). Both sheets share the same rows.
$wb = $excel.ActiveWorkbook

$updates = @{
    4  = 10819
    6  = 979
    8  = 1337
    9  = 8294
    10 = 39
    12 = 587
    13 = 219
    16 = 41
    18 = 30
    19 = 785
    24 = 1776
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
